# Automatische test-sync: 2025-06-18 12:00:10
# Adds the 11:30 "Afmelding nieuwsbrief" unsubscribe log entry to the
# "Logs" sheet and syncs the "Dashboard" category-count sheet to match.

$wb = $excel.ActiveWorkbook

# --- 1. Append the new log row to the "Logs" sheet -------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A9").Value = "Afmelding nieuwsbrief"
$logs.Range("B9").Value = "mailmind.test@zohomail.eu"
$logs.Range("C9").Value = "Graag afmelden voor de nieuwsbrief. Dank u."
$logs.Range("D9").Value = "Afmelding"
$logs.Range("F9").Value = "2025-06-18 11:30:12"
$logs.Range("G9").Value = "Nee"

# Extend the conditional formatting ranges so the new row is covered too,
# same as the existing rules (D2:D8 -> D2:D9, G2:G8 -> G2:G9).
$catFormats = $logs.Range("D2:D8").FormatConditions
for ($i = 1; $i -le $catFormats.Count; $i++) {
    $catFormats.Item($i).ModifyAppliesToRange($logs.Range("D2:D9"))
}

$answeredFormats = $logs.Range("G2:G8").FormatConditions
for ($i = 1; $i -le $answeredFormats.Count; $i++) {
    $answeredFormats.Item($i).ModifyAppliesToRange($logs.Range("G2:G9"))
}

# --- 2. Re-sync the "Dashboard" category counts -----------------------------
$dash = $wb.Worksheets.Item("Dashboard")

# "Afmelding" now has 2 occurrences and moves up to row 3; "Overig" (still 1)
# drops to row 4.
$dash.Range("A3").Value = "Afmelding"
$dash.Range("B3").Value = 2
$dash.Range("A4").Value = "Overig"
$dash.Range("B4").Value = 1
